# 20190716 done, 20190717 init

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 17 (20190716) - fill in the remaining checkmarks/crosses
$ws.Range("D17").Value = "×"
$ws.Range("G17").Value = "×"
$ws.Range("H17").Value = "×"
$ws.Range("M17").Value = "√"
$ws.Range("O17").Value = "×"

# Add new row 18 (20190717) - initialize with data recorded so far
$ws.Range("A18").Value = 20190717
$ws.Range("B18").Value = "√"
$ws.Range("C18").Value = "√"
$ws.Range("E18").Value = "√"
$ws.Range("F18").Value = "√"
$ws.Range("I18").Value = "√"
$ws.Range("J18").Value = "√"
$ws.Range("K18").Value = "√"
$ws.Range("O18").Value = "√"

# Update selection to reflect the new active cell
$ws.Range("O18").Select()
